# Fixed sample products cupcake typo
$wb = $excel.ActiveWorkbook

# --- Main sheet: fix the "Cup Cake Sample" product's name / meta fields ---
$main = $wb.Worksheets.Item("Main")
$main.Range("E6").Value = "Cupcake Sample"
$main.Range("O6").Value = "Vanilla Cupcake with Rich Frosting"
$main.Range("P6").Value = "Vanilla Cupcake with Rich Frosting"

# --- Categories sheet: update the product slug reference for the cupcake product ---
$cats = $wb.Worksheets.Item("Categories")
$cats.Range("A5").Value = "cupcake-sample"

# --- Make Categories the active sheet/tab, matching the saved selection state ---
$cats.Activate()
$cats.Range("A5").Select()
